$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Extend the "accelerating structure" lookup table (columns J:N) from
# row 65 down to row 113 to match rows already filled for 2:64 -
# J gets the running index (copied from column H), and K:N reproduce
# the same formula pattern already used in K3:N64.
# ------------------------------------------------------------------

# Read the already-computed H column values (D,E,F derived) for the
# rows we are about to fill in J - these are plain numbers, same as
# what J65:J113 should contain.
$hVals = $ws.Range("H65:H113").Value2

$ws.Range("J65:J113").Value2 = $hVals

$ws.Range("K65:K113").Formula = "=ROUND((-1 + SQRT(1 + 4 * N65)) / 2,0)"
$ws.Range("L65:L113").Formula = "=N65-K65*(K65+1)"
$ws.Range("M65:M113").Formula = "=MOD(J65,7)-`$B`$2"
$ws.Range("N65:N113").Formula = "=ROUNDDOWN(J65/(2*`$B`$2+1),0)"

# Column J (width bestFit) grew to fit 3-digit numbers.
$ws.Columns.Item(10).ColumnWidth = 4

# Update the view: scroll back to the top-left and select the newly
# filled block.
$ws.Range("A1").Select() | Out-Null
$ws.Range("K64:N113").Select() | Out-Null

Write-Host "done"
